# Chad-Itgb1 LR-pairs sheet: refresh with the new TPM-derived expression
# values. The sheet has one row per (sending cluster, target cluster) pair.
# Columns E:J describe the ligand (sending cluster) side, K:P describe the
# receptor (target cluster) side, and Q:T are derived edge weights -
# exactly the NATMI "lrc2p" computation:
#
#   I = G / sum(G over all sending clusters)      (ligand avg-expr specificity)
#   J = H / sum(H over all sending clusters)       (ligand total-expr specificity)
#   O = M / sum(M over all target clusters)        (receptor avg-expr specificity)
#   P = N / sum(N over all target clusters)        (receptor total-expr specificity)
#   Q = G * M   R = H * N                          (edge expression weights)
#   S = I * O   T = J * P                          (edge specificity weights)
#
# The commit updates the underlying per-cluster TPM figures (ligand average/
# total expression G/H, some ligand-expressing-cell counts E/F, and receptor
# average/total expression M/N), then every derived column is recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# --- New per-sending-cluster ligand figures (column A) -----------------
# Ligand-expressing cells (E), detection rate (F), average expr (G), total expr (H)
$senderNew = @{
    "ECs"              = @{ E = 2; F = 0.6666666666666666; G = 0.3484863333333333; H = 1.045459 }
    "FAPs"             = @{ E = 3; F = 1;                   G = 2.782805666666667; H = 8.348417 }
    "Inflammatory-Mac" = @{ E = 2; F = 0.6666666666666666; G = 0.5631773333333333; H = 1.689532 }
    "MuSCs"            = @{ E = 2; F = 0.6666666666666666; G = 0.2482816666666667; H = 0.744845 }
    "Resolving-Mac"    = @{ E = 2; F = 0.6666666666666666; G = 0.1494596666666667; H = 0.448379 }
}

# --- New per-target-cluster receptor figures (column D) -----------------
# Average expr (M), total expr (N)
$targetNew = @{
    "ECs"              = @{ M = 176.8550973333333;  N = 530.565292 }
    "FAPs"             = @{ M = 145.6413626666667;  N = 436.924088 }
    "Inflammatory-Mac" = @{ M = 84.02511333333334;  N = 252.07534 }
    "MuSCs"            = @{ M = 84.92877566666668;  N = 254.786327 }
    "Neutrophils"      = @{ M = 63.97102366666667;  N = 191.913071 }
    "Resolving-Mac"    = @{ M = 107.0290476666667;  N = 321.087143 }
}

# --- Write the updated base figures (E,F,G,H and M,N) first -------------
for ($r = 2; $r -le $lastRow; $r++) {
    $sender = $ws.Cells.Item($r, 1).Value2   # column A: Sending cluster
    $target = $ws.Cells.Item($r, 4).Value2   # column D: Target cluster

    $s = $senderNew[$sender]
    $t = $targetNew[$target]

    $ws.Cells.Item($r, 5).Value  = $s.E   # E
    $ws.Cells.Item($r, 6).Value  = $s.F   # F
    $ws.Cells.Item($r, 7).Value  = $s.G   # G
    $ws.Cells.Item($r, 8).Value  = $s.H   # H

    $ws.Cells.Item($r, 13).Value = $t.M   # M
    $ws.Cells.Item($r, 14).Value = $t.N   # N
}

# --- Totals across all (distinct) sending / target clusters -------------
$sumG = 0.0
$sumH = 0.0
foreach ($k in $senderNew.Keys) {
    $sumG += $senderNew[$k].G
    $sumH += $senderNew[$k].H
}

$sumM = 0.0
$sumN = 0.0
foreach ($k in $targetNew.Keys) {
    $sumM += $targetNew[$k].M
    $sumN += $targetNew[$k].N
}

# --- Recompute every derived column (I,J,O,P,Q,R,S,T) --------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    $h = $ws.Cells.Item($r, 8).Value2
    $m = $ws.Cells.Item($r, 13).Value2
    $n = $ws.Cells.Item($r, 14).Value2

    $i = $g / $sumG
    $j = $h / $sumH
    $o = $m / $sumM
    $p = $n / $sumN

    $ws.Cells.Item($r, 9).Value  = $i        # I
    $ws.Cells.Item($r, 10).Value = $j        # J
    $ws.Cells.Item($r, 15).Value = $o        # O
    $ws.Cells.Item($r, 16).Value = $p        # P

    $ws.Cells.Item($r, 17).Value = $g * $m   # Q
    $ws.Cells.Item($r, 18).Value = $h * $n   # R
    $ws.Cells.Item($r, 19).Value = $i * $o   # S
    $ws.Cells.Item($r, 20).Value = $j * $p   # T
}
